# Update header row (row 1) on the active sheet to reflect the new set of
# logged columns: added "E_generation_ep_steps", removed "M_state_size" and
# "M_action_size", and added four new columns for the environment / rule
# based model support (O_environment, O_NS_mult, O_EW_mult, O_phase_end_offset).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "time",
    "rew",
    "waitingTime",
    "E_num_train_rollouts",
    "E_rollout_length",
    "E_eval_freq",
    "E_eval_num_eps",
    "E_max_ep_steps",
    "E_generation_ep_steps",
    "E_test_num_eps",
    "A_gae_tau",
    "A_entropy_weight",
    "A_minibatch_size",
    "A_optimization_epochs",
    "A_ppo_ratio_clip",
    "A_discount",
    "A_learning_rate",
    "A_clip_grads",
    "A_gradient_clip",
    "A_value_loss_coef",
    "O_num_agents",
    "O_rule_set",
    "O_rule_set_params",
    "O_environment",
    "O_NS_mult",
    "O_EW_mult",
    "O_phase_end_offset"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# Make sure the newly added header cells (Z1, AA1) carry the same style as
# the rest of the header row (bold font + border), matching the formatting
# already applied to the other header cells.
$headerStyleSource = $ws.Range("Y1")
$newHeaderCells = $ws.Range("Z1:AA1")
$headerStyleSource.Copy()
$newHeaderCells.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
